$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Project Planner")

# --- Fill in "Actual Start" (col E) / "Actual Duration" (col F) progress values ---
$ws.Range("E5").Value  = 1
$ws.Range("F5").Value  = 2

$ws.Range("E6").Value  = 1
$ws.Range("F6").Value  = 2

$ws.Range("F7").Value  = 3

$ws.Range("E8").Value  = 2
$ws.Range("F8").Value  = 2

$ws.Range("E9").Value  = 2
$ws.Range("F9").Value  = 3

$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 3

$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 3

$ws.Range("E16").Value = 4
$ws.Range("F16").Value = 2

$ws.Range("E18").Value = 4

$ws.Range("E21").Value = 5

$ws.Range("E24").Value = 5

# --- Update the view: scroll/zoom out a bit and move the selection ---
$ws.Activate()
$excel.ActiveWindow.Zoom = 90
$ws.Range("F11").Select()
